$wb = $excel.ActiveWorkbook

Write-Host "Updating sheet ALC..."
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 495.81818
$ws.Range("J17").Value = 495.81818
$ws.Range("L17").Value = 1487.45454
$ws.Range("N17").Value = -1823.45454
$ws.Range("H62").Value = 2777.3333
$ws.Range("I62").Value = 1463.3334
$ws.Range("J62").Value = 3434.3333
$ws.Range("K62").Value = 1463.3334
$ws.Range("L62").Value = 3434.3333
$ws.Range("M62").Value = -839.3334
$ws.Range("N62").Value = -4682.3333
$ws.Range("H65").Value = 2777.3333
$ws.Range("I65").Value = 1463.3334
$ws.Range("J65").Value = 3434.3333
$ws.Range("K65").Value = 7316.666999999999
$ws.Range("L65").Value = 17171.6665
$ws.Range("M65").Value = -4196.666999999999
$ws.Range("N65").Value = -23411.6665
$ws.Range("H92").Value = 612.86957
$ws.Range("I92").Value = 476.92856
$ws.Range("J92").Value = 824.3333
$ws.Range("K92").Value = 476.92856
$ws.Range("L92").Value = 824.3333
$ws.Range("M92").Value = 771.0714399999999
$ws.Range("N92").Value = -3320.3333
$ws.Range("H96").Value = 347.5
$ws.Range("I96").Value = 354.7143
$ws.Range("J96").Value = 337.4
$ws.Range("K96").Value = 1064.1429
$ws.Range("L96").Value = 1012.2
$ws.Range("M96").Value = 308.8571000000002
$ws.Range("N96").Value = -3758.2
$ws.Range("H97").Value = 1608.5714
$ws.Range("I97").Value = 740
$ws.Range("J97").Value = 1956
$ws.Range("K97").Value = 2220
$ws.Range("L97").Value = 5868
$ws.Range("M97").Value = -1724
$ws.Range("N97").Value = -6860
$ws.Range("H113").Value = 3265.3
$ws.Range("I113").Value = 3053.3333
$ws.Range("J113").Value = 3901.2
$ws.Range("K113").Value = 3053.3333
$ws.Range("L113").Value = 3901.2
$ws.Range("M113").Value = 200.6667000000002
$ws.Range("N113").Value = -10409.2

Write-Host "Updating sheet ARM..."
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1528.15
$ws.Range("I32").Value = 1095.625
$ws.Range("J32").Value = 2640.3572
$ws.Range("K32").Value = 1095.625
$ws.Range("L32").Value = 2640.3572
$ws.Range("M32").Value = -808.625
$ws.Range("N32").Value = -3214.3572
$ws.Range("H61").Value = 111334460
$ws.Range("I61").Value = 200200640
$ws.Range("J61").Value = 251750
$ws.Range("K61").Value = 200200640
$ws.Range("L61").Value = 251750
$ws.Range("M61").Value = -200200428
$ws.Range("N61").Value = -252174
$ws.Range("H97").Value = 2718214.2
$ws.Range("I97").Value = 4808517.5
$ws.Range("J97").Value = 820
$ws.Range("K97").Value = 4808517.5
$ws.Range("L97").Value = 820
$ws.Range("M97").Value = -4808021.5
$ws.Range("N97").Value = -1812
$ws.Range("H102").Value = 6217342
$ws.Range("I102").Value = 7149493.5
$ws.Range("J102").Value = 3000
$ws.Range("K102").Value = 7149493.5
$ws.Range("L102").Value = 3000
$ws.Range("M102").Value = -7147871.5
$ws.Range("N102").Value = -6244
$ws.Range("H135").Value = 32662.691
$ws.Range("J135").Value = 32662.691
$ws.Range("L135").Value = 32662.691
$ws.Range("N135").Value = -42802.691
$ws.Range("H136").Value = 111334460
$ws.Range("I136").Value = 200200640
$ws.Range("J136").Value = 251750
$ws.Range("K136").Value = 600601920
$ws.Range("L136").Value = 755250
$ws.Range("M136").Value = -600599370
$ws.Range("N136").Value = -760350

Write-Host "Updating sheet BSM..."
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H19").Value = 0
$ws.Range("J19").Value = 0
$ws.Range("L19").Value = 0
$ws.Range("N19").ClearContents()
$ws.Range("H94").Value = 487
$ws.Range("I94").Value = 357.18182
$ws.Range("J94").Value = 725
$ws.Range("K94").Value = 357.18182
$ws.Range("L94").Value = 725
$ws.Range("M94").Value = 93.81817999999998
$ws.Range("N94").Value = -1627
$ws.Range("H99").Value = 1500
$ws.Range("J99").Value = 1500
$ws.Range("L99").Value = 1500
$ws.Range("N99").Value = -4496
$ws.Range("H105").Value = 35716144
$ws.Range("I105").Value = 62501664
$ws.Range("J105").Value = 2116.6667
$ws.Range("K105").Value = 62501664
$ws.Range("L105").Value = 2116.6667
$ws.Range("M105").Value = -62499917
$ws.Range("N105").Value = -5610.6667

Write-Host "Updating sheet CRP..."
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4131.1787
$ws.Range("I31").Value = 1492.9445
$ws.Range("J31").Value = 8880
$ws.Range("K31").Value = 1492.9445
$ws.Range("L31").Value = 8880
$ws.Range("M31").Value = -1197.9445
$ws.Range("N31").Value = -9470
$ws.Range("H34").Value = 4131.1787
$ws.Range("I34").Value = 1492.9445
$ws.Range("J34").Value = 8880
$ws.Range("K34").Value = 1492.9445
$ws.Range("L34").Value = 8880
$ws.Range("M34").Value = -1290.9445
$ws.Range("N34").Value = -9284
$ws.Range("H58").Value = 19609448
$ws.Range("I58").Value = 23810932
$ws.Range("J58").Value = 2520.889
$ws.Range("K58").Value = 23810932
$ws.Range("L58").Value = 2520.889
$ws.Range("M58").Value = -23810729
$ws.Range("N58").Value = -2926.889
$ws.Range("H74").Value = 25313.842
$ws.Range("J74").Value = 25313.842
$ws.Range("L74").Value = 25313.842
$ws.Range("N74").Value = -27061.842
$ws.Range("H75").Value = 49800
$ws.Range("J75").Value = 49800
$ws.Range("L75").Value = 49800
$ws.Range("N75").Value = -51796
$ws.Range("H77").Value = 25313.842
$ws.Range("J77").Value = 25313.842
$ws.Range("L77").Value = 75941.526
$ws.Range("N77").Value = -84677.526
$ws.Range("H78").Value = 49800
$ws.Range("J78").Value = 49800
$ws.Range("L78").Value = 149400
$ws.Range("N78").Value = -159384
$ws.Range("H107").Value = 298.75
$ws.Range("I107").Value = 309.0476
$ws.Range("J107").Value = 226.66667
$ws.Range("K107").Value = 309.0476
$ws.Range("L107").Value = 226.66667
$ws.Range("M107").Value = 1610.9524
$ws.Range("N107").Value = -4066.66667
$ws.Range("H136").Value = 19609448
$ws.Range("I136").Value = 23810932
$ws.Range("J136").Value = 2520.889
$ws.Range("K136").Value = 71432796
$ws.Range("L136").Value = 7562.667
$ws.Range("M136").Value = -71430246
$ws.Range("N136").Value = -12662.667

Write-Host "Updating sheet CUL..."
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 437.25
$ws.Range("I98").Value = 83
$ws.Range("J98").Value = 1500
$ws.Range("K98").Value = 249
$ws.Range("L98").Value = 4500
$ws.Range("M98").Value = 1249
$ws.Range("N98").Value = -7496
$ws.Range("H131").Value = 905.322
$ws.Range("J131").Value = 1025.2766
$ws.Range("L131").Value = 3075.8298
$ws.Range("N131").Value = -13155.8298
$ws.Range("H136").Value = 2338.0715
$ws.Range("I136").Value = 1375.7142
$ws.Range("J136").Value = 3300.4285
$ws.Range("K136").Value = 4127.142599999999
$ws.Range("L136").Value = 9901.2855
$ws.Range("M136").Value = 972.8574000000008
$ws.Range("N136").Value = -20101.2855
$ws.Range("H140").Value = 2334.4905
$ws.Range("I140").Value = 2316.318
$ws.Range("J140").Value = 2347.3872
$ws.Range("K140").Value = 6948.954000000001
$ws.Range("L140").Value = 7042.1616
$ws.Range("M140").Value = -1768.954000000001
$ws.Range("N140").Value = -17402.1616

Write-Host "Updating sheet GSM..."
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H93").Value = 0
$ws.Range("I93").Value = 0
$ws.Range("K93").Value = 0
$ws.Range("M93").ClearContents()
$ws.Range("H113").Value = 1376.4375
$ws.Range("I113").Value = 1128.5714
$ws.Range("J113").Value = 1569.2222
$ws.Range("K113").Value = 1128.5714
$ws.Range("L113").Value = 1569.2222
$ws.Range("M113").Value = 1041.4286
$ws.Range("N113").Value = -5909.2222
$ws.Range("H132").Value = 112825.39
$ws.Range("I132").Value = 64178.562
$ws.Range("J132").Value = 502000
$ws.Range("K132").Value = 192535.686
$ws.Range("L132").Value = 1506000
$ws.Range("M132").Value = -190005.686
$ws.Range("N132").Value = -1511060
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()
$ws.Range("H137").Value = 0
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("M137").ClearContents()
$ws.Range("N137").ClearContents()
$ws.Range("H138").Value = 100000
$ws.Range("J138").Value = 100000
$ws.Range("L138").Value = 100000
$ws.Range("N138").Value = -110280
$ws.Range("H139").Value = 100000
$ws.Range("J139").Value = 100000
$ws.Range("L139").Value = 100000
$ws.Range("N139").Value = -4496

Write-Host "Updating sheet LTW..."
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2886.125
$ws.Range("I61").Value = 2897.5454
$ws.Range("K61").Value = 2897.5454
$ws.Range("M61").Value = -2695.5454
$ws.Range("H93").Value = 823.4583
$ws.Range("I93").Value = 756.5
$ws.Range("J93").Value = 957.375
$ws.Range("K93").Value = 756.5
$ws.Range("L93").Value = 957.375
$ws.Range("M93").Value = 491.5
$ws.Range("N93").Value = -3453.375
$ws.Range("H113").Value = 2886.125
$ws.Range("I113").Value = 2897.5454
$ws.Range("K113").Value = 2897.5454
$ws.Range("M113").Value = -727.5454
$ws.Range("H136").Value = 90898.96000000001
$ws.Range("I136").Value = 70213.75
$ws.Range("J136").Value = 132269.38
$ws.Range("K136").Value = 210641.25
$ws.Range("L136").Value = 396808.14
$ws.Range("M136").Value = -208091.25
$ws.Range("N136").Value = -401908.14

Write-Host "Updating sheet WVR..."
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 3250
$ws.Range("J5").Value = 3250
$ws.Range("L5").Value = 3250
$ws.Range("N5").Value = -3474
$ws.Range("H96").Value = 1692.3334
$ws.Range("I96").Value = 1087.5
$ws.Range("J96").Value = 2902
$ws.Range("K96").Value = 1087.5
$ws.Range("L96").Value = 2902
$ws.Range("M96").Value = 285.5
$ws.Range("N96").Value = -5648
$ws.Range("H100").Value = 67909.8
$ws.Range("I100").Value = 100656.4
$ws.Range("J100").Value = 51536.5
$ws.Range("K100").Value = 201312.8
$ws.Range("L100").Value = 103073
$ws.Range("M100").Value = -200771.8
$ws.Range("N100").Value = -104155

Write-Host "Done updating Hades_Profits sheets."
